# Update "Inscritos" (E), "Pagos" (F) and "Inscricoes homologadas" (H) figures
# in the Inscricoes summary table to reflect newly processed registrations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: E=5 (Inscritos), F=6 (Pagos), G=7 (Isencoes deferidas), H=8 (Inscricoes homologadas)

$updates = @(
    @{ Row = 5;  E = 121; F = 82;  H = 82  },
    @{ Row = 6;  F = 28;  H = 28  },
    @{ Row = 10; E = 423; F = 211; H = 211 },
    @{ Row = 11; E = 285; F = 159; H = 159 },
    @{ Row = 12; E = 410; F = 223; H = 223 },
    @{ Row = 14; F = 52;  H = 52  },
    @{ Row = 15; F = 54;  H = 54  },
    @{ Row = 16; F = 90;  H = 90  },
    @{ Row = 17; F = 41;  H = 41  },
    @{ Row = 18; E = 49;  F = 23;  H = 23  },
    @{ Row = 23; F = 77;  H = 77  },
    @{ Row = 24; E = 181 },
    @{ Row = 25; E = 220; F = 102; H = 102 },
    @{ Row = 26; F = 77;  H = 77  },
    @{ Row = 27; E = 281; F = 134; H = 134 },
    @{ Row = 28; F = 57;  H = 57  },
    @{ Row = 30; E = 184; F = 105; H = 105 },
    @{ Row = 32; F = 97;  H = 97  },
    @{ Row = 33; F = 128; H = 128 },
    @{ Row = 34; F = 113; H = 113 },
    @{ Row = 35; E = 121 },
    @{ Row = 39; E = 165 },
    @{ Row = 40; E = 226; F = 105; H = 105 },
    @{ Row = 42; F = 164; H = 164 },
    @{ Row = 43; F = 54;  H = 54  },
    @{ Row = 44; E = 266; F = 128; H = 128 },
    @{ Row = 45; E = 120 },
    @{ Row = 46; E = 268; F = 147; H = 147 },
    @{ Row = 47; F = 188; H = 188 },
    @{ Row = 49; E = 255; F = 107; H = 107 },
    @{ Row = 50; E = 218; F = 95;  H = 95  },
    @{ Row = 51; E = 203; F = 84;  H = 84  }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
    if ($u.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $u.F }
    if ($u.ContainsKey('H')) { $ws.Cells.Item($r, 8).Value = $u.H }
}
